# Clean up footnote markers ("[1]".."[5]") and embedded line breaks in the
# vaccine/brand/description labels across every sheet of the workbook.
#
#   "DTaP [1]"                                   -> "DTaP "
#   "Hepatitis B [5]\nPediatric/Adolescent"       -> "Hepatitis B  Pediatric/Adolescent"
#   "Recombivax\nHB"                              -> "Recombivax HB"
#   "Tetanus and Diphtheria Toxoids[3]"           -> "Tetanus and Diphtheria Toxoids"
#   ... etc.

$wb = $excel.ActiveWorkbook

$newline = [char]10

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Remove the footnote markers "[1]" through "[5]" wherever they appear
    # (LookAt:=2 is xlPart, so it matches the bracketed text inside a cell).
    foreach ($n in 1..5) {
        $marker = "[" + $n + "]"
        [void]$ws.Cells.Replace($marker, "", 2, 1, $false, $false)
    }

    # Collapse embedded line breaks within a cell into a single space so
    # multi-line labels become one line of text.
    [void]$ws.Cells.Replace($newline, " ", 2, 1, $false, $false)
}
